# Auto-generated edit script: applies leve-profit recalculations
# across ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2650.1667
$ws.Range("I51").Value = 2475.25
$ws.Range("K51").Value = 2475.25
$ws.Range("M51").Value = -1991.25

$ws.Range("H125").Value = 4512.222
$ws.Range("I125").Value = 5332.4
$ws.Range("J125").Value = 3487
$ws.Range("K125").Value = 47991.6
$ws.Range("L125").Value = 31383
$ws.Range("M125").Value = -45531.6
$ws.Range("N125").Value = -36303

$ws.Range("H132").Value = 1764.1538
$ws.Range("I132").Value = 1802.72
$ws.Range("K132").Value = 5408.16
$ws.Range("M132").Value = -2878.16

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9763.5
$ws.Range("I45").Value = 10540.728
$ws.Range("J45").Value = 1214
$ws.Range("K45").Value = 10540.728
$ws.Range("L45").Value = 1214
$ws.Range("M45").Value = -10163.728
$ws.Range("N45").Value = -1968

$ws.Range("H61").Value = 4806.1113
$ws.Range("I61").Value = 7489
$ws.Range("J61").Value = 3349.6858
$ws.Range("K61").Value = 7489
$ws.Range("L61").Value = 3349.6858
$ws.Range("M61").Value = -7277
$ws.Range("N61").Value = -3773.6858

$ws.Range("H63").Value = 62502572
$ws.Range("I63").Value = 90910290
$ws.Range("J63").Value = 5601.2
$ws.Range("K63").Value = 90910290
$ws.Range("L63").Value = 5601.2
$ws.Range("M63").Value = -90909604
$ws.Range("N63").Value = -6973.2

$ws.Range("H66").Value = 62502572
$ws.Range("I66").Value = 90910290
$ws.Range("J66").Value = 5601.2
$ws.Range("K66").Value = 454551450
$ws.Range("L66").Value = 28006
$ws.Range("M66").Value = -454548018
$ws.Range("N66").Value = -34870

$ws.Range("H122").Value = 1604861.5
$ws.Range("I122").Value = 1833727.5
$ws.Range("K122").Value = 5501182.5
$ws.Range("M122").Value = -5498732.5

$ws.Range("H132").Value = 2521.3274
$ws.Range("I132").Value = 1298.6154
$ws.Range("K132").Value = 3895.8462
$ws.Range("M132").Value = -1365.8462

$ws.Range("H136").Value = 4806.1113
$ws.Range("I136").Value = 7489
$ws.Range("J136").Value = 3349.6858
$ws.Range("K136").Value = 22467
$ws.Range("L136").Value = 10049.0574
$ws.Range("M136").Value = -19917
$ws.Range("N136").Value = -15149.0574

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 76944264
$ws.Range("I20").Value = 250001700
$ws.Range("J20").Value = 29846.555
$ws.Range("K20").Value = 250001700
$ws.Range("L20").Value = 29846.555
$ws.Range("M20").Value = -250001453
$ws.Range("N20").Value = -30340.555

$ws.Range("H94").Value = 1381.12
$ws.Range("I94").Value = 886.9375
$ws.Range("J94").Value = 2259.6667
$ws.Range("K94").Value = 886.9375
$ws.Range("L94").Value = 2259.6667
$ws.Range("M94").Value = -435.9375
$ws.Range("N94").Value = -3161.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2901.9285
$ws.Range("I31").Value = 1865.6842
$ws.Range("J31").Value = 4132.4688
$ws.Range("K31").Value = 1865.6842
$ws.Range("L31").Value = 4132.4688
$ws.Range("M31").Value = -1570.6842
$ws.Range("N31").Value = -4722.4688

$ws.Range("H34").Value = 2901.9285
$ws.Range("I34").Value = 1865.6842
$ws.Range("J34").Value = 4132.4688
$ws.Range("K34").Value = 1865.6842
$ws.Range("L34").Value = 4132.4688
$ws.Range("M34").Value = -1663.6842
$ws.Range("N34").Value = -4536.4688

$ws.Range("H86").Value = 125002470
$ws.Range("I86").Value = 250002480
$ws.Range("J86").Value = 2475
$ws.Range("K86").Value = 250002480
$ws.Range("L86").Value = 2475
$ws.Range("M86").Value = -250001357
$ws.Range("N86").Value = -4721

$ws.Range("H89").Value = 125002470
$ws.Range("I89").Value = 250002480
$ws.Range("J89").Value = 2475
$ws.Range("K89").Value = 1250012400
$ws.Range("L89").Value = 12375
$ws.Range("M89").Value = -1250006784
$ws.Range("N89").Value = -23607

$ws.Range("H99").Value = 8940641
$ws.Range("I99").Value = 14997.272
$ws.Range("J99").Value = 41668000
$ws.Range("K99").Value = 14997.272
$ws.Range("L99").Value = 41668000
$ws.Range("M99").Value = -13499.272
$ws.Range("N99").Value = -41670996

$ws.Range("H126").Value = 8940641
$ws.Range("I126").Value = 14997.272
$ws.Range("J126").Value = 41668000
$ws.Range("K126").Value = 44991.81600000001
$ws.Range("L126").Value = 125004000
$ws.Range("M126").Value = -42521.81600000001
$ws.Range("N126").Value = -125008940

$ws.Range("H134").Value = 3419.3225
$ws.Range("I134").Value = 4210.263
$ws.Range("J134").Value = 2167
$ws.Range("K134").Value = 12630.789
$ws.Range("L134").Value = 6501
$ws.Range("M134").Value = -10095.789
$ws.Range("N134").Value = -11571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H122").Value = 21371214
$ws.Range("I122").Value = 2781098.2
$ws.Range("J122").Value = 83338264
$ws.Range("K122").Value = 8343294.600000001
$ws.Range("L122").Value = 250014792
$ws.Range("M122").Value = -8340844.600000001
$ws.Range("N122").Value = -250019692

$ws.Range("H123").Value = 18453.97
$ws.Range("J123").Value = 18740.455
$ws.Range("L123").Value = 18740.455
$ws.Range("N123").Value = -23640.455

$ws.Range("H126").Value = 7606.222
$ws.Range("I126").Value = 8467.467000000001
$ws.Range("J126").Value = 3300
$ws.Range("K126").Value = 25402.401
$ws.Range("L126").Value = 9900
$ws.Range("M126").Value = -22932.401
$ws.Range("N126").Value = -14840

$ws.Range("H132").Value = 29406.37
$ws.Range("I132").Value = 70713.47
$ws.Range("J132").Value = 2466.9565
$ws.Range("K132").Value = 212140.41
$ws.Range("L132").Value = 7400.869499999999
$ws.Range("M132").Value = -209610.41
$ws.Range("N132").Value = -12460.8695

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 50001890
$ws.Range("I68").Value = 1952.4546
$ws.Range("J68").Value = 111112920
$ws.Range("K68").Value = 1952.4546
$ws.Range("L68").Value = 111112920
$ws.Range("M68").Value = -1203.4546
$ws.Range("N68").Value = -111114418

$ws.Range("H71").Value = 50001890
$ws.Range("I71").Value = 1952.4546
$ws.Range("J71").Value = 111112920
$ws.Range("K71").Value = 9762.273000000001
$ws.Range("L71").Value = 555564600
$ws.Range("M71").Value = -6018.273000000001
$ws.Range("N71").Value = -555572088

$ws.Range("H93").Value = 1650.75
$ws.Range("I93").Value = 1701.5
$ws.Range("K93").Value = 1701.5
$ws.Range("M93").Value = -453.5

$ws.Range("H132").Value = 15880553
$ws.Range("I132").Value = 22231420
$ws.Range("J132").Value = 3384.6667
$ws.Range("K132").Value = 66694260
$ws.Range("L132").Value = 10154.0001
$ws.Range("M132").Value = -66691730
$ws.Range("N132").Value = -15214.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2162.5
$ws.Range("I122").Value = 2118.75
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 6356.25
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -3906.25
$ws.Range("N122").Value = -11650

$ws.Range("H136").Value = 3408.2727
$ws.Range("I136").Value = 4359.385
$ws.Range("J136").Value = 2034.4445
$ws.Range("K136").Value = 13078.155
$ws.Range("L136").Value = 6103.333500000001
$ws.Range("M136").Value = -10528.155
$ws.Range("N136").Value = -11203.3335
